$d = $word.ActiveDocument
$errs = $d.SpellingErrors
Write-Output $errs.Count
